# This workbook is a 2-column Q&A dataset ("surprised" sentiment corpus):
# column A holds a question/title string per row, column B always holds "1".
# The commit removes 81 off-topic rows (mostly rows whose text references the
# unrelated TV show "기막힌외출") that had slipped into the corpus.
#
# Row numbers below are the ORIGINAL (pre-edit) 1-based worksheet row numbers
# of the entries being removed, sorted ascending. We delete from the bottom
# of the sheet upward so that removing one row never shifts the row number
# of another row still queued for deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(
    16, 30, 31, 53, 75, 190, 278, 311, 320, 349, 395, 477, 480,
    497, 668, 725, 752, 764, 817, 855, 871, 875, 1085, 1194, 1218, 1364,
    1527, 1550, 1588, 1646, 1688, 1718, 1767, 1790, 1819, 1851, 1894, 1939, 1989,
    2083, 2112, 2118, 2199, 2210, 2233, 2246, 2250, 2263, 2330, 2368, 2530, 2704,
    2728, 2804, 2850, 2929, 3260, 3305, 3308, 3349, 3395, 3492, 3548, 3569, 3603,
    3662, 3688, 3709, 3755, 3821, 3851, 3905, 3936, 3952, 4053, 4123, 4137, 4193,
    4234, 4247, 4424
)

$sortedDescending = $rowsToDelete | Sort-Object -Descending

foreach ($rowNum in $sortedDescending) {
    $ws.Rows($rowNum).Delete()
}

# Restore the workbook's on-screen selection/scroll target recorded in the
# saved file after the edit (final active cell sits in column B near the
# bottom of the now-shorter sheet).
$ws.Range("B4164").Select()
